# gruposPlantilla.xlsx: replace the single "ASIGNATURAS" title row with a
# full header row that includes a new "Código Asignatura" column (first)
# and a renamed "Curso" column (was "Horario", now last), and move the
# small underline-styled marker cell accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a brand-new column at A. This shifts the old A:E range to B:F,
#    which means:
#      - the merged title cell A1:E1 becomes B1:F1
#      - the old header row (row 2) becomes the B2:F2 range
#      - the little marker cell (old D11) becomes E11
#    and, importantly, the original column width (20.77734375) carried by
#    columns A:D is preserved unmodified on the new B:E columns.
$ws.Columns.Item(1).Insert()

# 2. The old title cell was merged across B1:F1 - split it back apart.
$ws.Range("B1:F1").UnMerge()

# 3. Drop the old title row (row 1, "ASIGNATURAS") entirely. Everything
#    below shifts up, so the former header row (Tipo/Grupo/Cuatrimestre/
#    Acreditacion/Horario) becomes row 1, and the marker cell becomes E10.
$ws.Rows.Item(1).Delete()

# 4. Column A is currently blank with default formatting. Give it the same
#    look as the rest of the header row (fill/font) by copying the format
#    from B1.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# 5. Fix up the text: rename the trailing header and add the new first
#    column header (set F1 first so the new shared string for A1 is
#    appended last, matching the source order).
$ws.Range("F1").Value2 = "Curso"
$ws.Range("A1").Value2 = "Código Asignatura"

# 6. Column A gets its own (slightly wider) width; B:E already inherited
#    the correct width from step 1 and don't need to be touched.
$ws.Columns.Item(1).ColumnWidth = 22.33203125

# 7. Restore the last active selection.
$ws.Range("B12").Select()
